# blacklisted_words.xlsx - keep only the French ("fra") rows, clearing the
# English and Arabic rows but leaving their row/column formatting in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two surviving data rows (formerly rows 6 & 7, the "fra" entries) move
# up to rows 2 & 3 and lose their special cell styling (plain/default style),
# with is_active now stored as a real boolean instead of the text "TRUE".
$ws.Range("A2:D3").Style = "Normal"

$ws.Range("A2").Value = "fra"
$ws.Range("B2").Value = "Merde"
$ws.Range("C2").Value = "Mot sur la liste noire"
$ws.Range("D2").Value = $true

$ws.Range("A3").Value = "fra"
$ws.Range("B3").Value = "bon sang"
$ws.Range("C3").Value = "Mot sur la liste noire"
$ws.Range("D3").Value = $true

# The remaining old rows (previously holding the English/eng and
# Arabic/ara blacklisted words) are emptied out, but keep their existing
# cell formatting (borders/fill/font) intact.
$ws.Range("A4:D10").ClearContents()

# Match the workbook's final on-screen selection.
$ws.Range("A2:D3").Select()
